# Auto-generated COM-interop script
# Applies updated percentages to 'iModulon' sheet (col B only)
# and a full relabel/resort + updated percentages to 'Subsystem' sheet
$wb = $excel.ActiveWorkbook

# --- iModulon sheet: update BAR_Set_percentage (column B) values only ---
$ws1 = $wb.Worksheets.Item("iModulon")
$iModulonB = @(
  1.147227533460803,
  1.720841300191205,
  8.030592734225621,
  7.648183556405354,
  0,
  2.103250478011472,
  1.720841300191205,
  0.9560229445506693,
  15.48757170172084,
  0.7648183556405354,
  0.5736137667304015,
  0.9560229445506693,
  3.441682600382409
)
$r = 2
foreach ($val in $iModulonB) {
    $ws1.Cells.Item($r, 2).Value2 = $val
    $r++
}

# --- Subsystem sheet: full relabel/resort with updated percentages ---
$ws2 = $wb.Worksheets.Item("Subsystem")
$subsystemRows = @(
  @('Alanine, aspartate and glutamate metabolism', 1.543739279588336, 2.173913043478261),
  @('Aminosugars metabolism', 0.8576329331046313, 0.7246376811594203),
  @('Arginine and proline metabolism', 2.572898799313894, 10.14492753623188),
  @('Biomass and maintenance functions', 0.1715265866209262, 0),
  @('Biotin metabolism', 0, 2.173913043478261),
  @('C5-Branched dibasic acid metabolism', 0.5145797598627788, 0),
  @('Carbon fixation', 1.029159519725558, 0),
  @('Carotenoid Biosynthesis', 2.229845626072041, 0),
  @('Citrate cycle (TCA cycle)', 1.029159519725558, 0),
  @('Cyanophycin metabolism', 0.1715265866209262, 0),
  @('Extracellular exchange', 3.259005145797599, 0),
  @('Fatty acid biosynthesis', 17.32418524871355, 0),
  @('Folate biosynthesis', 2.229845626072041, 1.449275362318841),
  @('Fructose and mannose metabolism', 0, 5.797101449275362),
  @('Galactolipids metabolism', 4.631217838765009, 0),
  @('Glutamate metabolism', 1.200686106346484, 0),
  @('Glutathione metabolism', 0, 1.449275362318841),
  @('Glycerolipid metabolism', 0.3430531732418525, 1.449275362318841),
  @('Glycolysis/Gluconeogenesis', 3.430531732418525, 4.347826086956522),
  @('Glyoxylate and dicarboxylate metabolism', 2.058319039451115, 1.449275362318841),
  @('Histidine metabolism', 1.543739279588336, 1.449275362318841),
  @('Hydrogen production', 0, 0),
  @('Inositol phosphate metabolism', 0, 1.449275362318841),
  @('Intracellular demand', 0.1715265866209262, 0),
  @('Intracellular source/sink', 0.1715265866209262, 0),
  @('Lipopolysaccharide biosynthesis', 0.8576329331046313, 0),
  @('Lysine metabolism', 1.543739279588336, 1.449275362318841),
  @('Nicotinate and nicotinamide metabolism', 0.8576329331046313, 1.449275362318841),
  @('Nitrogen metabolism', 1.715265866209263, 0.7246376811594203),
  @('Nucleotide sugars metabolism', 0, 2.898550724637681),
  @('Others', 0.6861063464837049, 1.449275362318841),
  @('Oxidative phosphorylation', 0.8576329331046313, 1.449275362318841),
  @('PHB byosynthesis', 0, 2.173913043478261),
  @('Pantothenate and CoA biosynthesis', 1.543739279588336, 0),
  @('Pentose phosphate pathway', 1.37221269296741, 0.7246376811594203),
  @('Peptidoglycan biosynthesis', 1.37221269296741, 0.7246376811594203),
  @('Phenylalanine tyrosine and tryptophan biosynthesis', 3.259005145797599, 5.797101449275362),
  @('Photosynthesis', 0.8576329331046313, 0),
  @('Porphyrin and chlorophyll metabolism', 7.032590051457976, 5.072463768115942),
  @('Purine metabolism', 4.459691252144083, 5.797101449275362),
  @('Pyrimidine metabolism', 3.430531732418525, 0.7246376811594203),
  @('Pyruvate metabolism', 1.715265866209263, 0.7246376811594203),
  @('Riboflavin metabolism', 1.715265866209263, 0),
  @('Starch and sucrose metabolism', 0.6861063464837049, 3.623188405797102),
  @('Steroid biosynthesis', 0, 2.173913043478261),
  @('Sterol biosynthesis', 1.029159519725558, 0),
  @('Sulfolipid Biosynthesis', 1.715265866209263, 0),
  @('Sulfur Cysteine and methionine metabolism', 2.401372212692967, 3.623188405797102),
  @('Terpenoid backbone biosynthesis', 1.543739279588336, 0),
  @('Thiamine metabolism', 0, 2.898550724637681),
  @('Transport', 7.204116638078903, 15.21739130434783),
  @('Ubiquinone and other pterpenoids biosynthesis', 2.401372212692967, 0),
  @('Urea cycle and metabolism of amino groups', 1.200686106346484, 0.7246376811594203),
  @('Valine leucine and isoleucine biosynthesis', 2.058319039451115, 2.898550724637681),
  @('Vitamin B6 metabolism', 0, 3.623188405797102)
)
$r = 2
foreach ($row in $subsystemRows) {
    $ws2.Cells.Item($r, 1).Value2 = $row[0]
    $ws2.Cells.Item($r, 2).Value2 = $row[1]
    $ws2.Cells.Item($r, 3).Value2 = $row[2]
    $r++
}

# New rows (53-56) need the same bold/border/center-top style as the
# existing label column cells (style copied from A2, which already
# carries the shared header/label style).
$ws2.Range("A2").Copy()
$ws2.Range("A53:A56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Done: updated iModulon and Subsystem sheets"
